$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Split the old combined subtitle text ("Empresa: ... Tema: Espacio")
#    into two separate strings: the subtitle keeps everything up to and
#    including "estratégico. " and a new row gets "Tema: Espacio".
$ws.Range("A2").Value = "Empresa: YYYYYYYYYY. Programa: Pensamiento estratégico. "

# 2. Make room for the new title-style row right below the subtitle.
#    Insert two blank rows at row 3, then drop the old (now shifted)
#    blank spacer row so the header/data rows land two rows further
#    down than before, with an empty row 4 gap in between.
$ws.Rows("3:4").Insert()
$ws.Rows("5").Delete()

# 3. Fill in the new row 3 with its own merged, title-styled text.
$ws.Range("A3").Value = "Tema: Espacio"
$ws.Range("A3:N3").Merge() | Out-Null
$ws.Range("A3:N3").Font.Name = "Arial"
$ws.Range("A3:N3").Font.Bold = $true
$ws.Range("A3:N3").Font.Size = 16
$ws.Range("A3:N3").HorizontalAlignment = -4108
$ws.Range("A3:N3").RowHeight = 20.25

# 4. Row 4 is just a spacer gap now - clear it out completely so it
#    doesn't carry any leftover formatting.
$ws.Rows("4").Clear() | Out-Null

# 5. Leftover touched cell far below the table (matches the stray
#    formatted cell seen at F28 in the final workbook).
$ws.Range("F28").IndentLevel = 0

# 6. Leave the selection where it ended up in the edited workbook.
$ws.Range("G12").Select() | Out-Null
